$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Add a brand new data row 4 ---
$ws.Cells.Item(4, 1).Value  = 42633.886377314811      # A4 - Date
$ws.Cells.Item(4, 1).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(4, 2).Value  = 3                        # B4 - ScoreFinal
$ws.Cells.Item(4, 3).Value  = "Neutral"                 # C4 - Verdict
$ws.Cells.Item(4, 4).Value  = -2                       # D4 - totalSentiment
$ws.Cells.Item(4, 5).Value  = 18232                    # E4 - wordCount
$ws.Cells.Item(4, 6).Value  = 1956                     # F4 - sentenceCount
$ws.Cells.Item(4, 7).Value  = 50                       # G4 - posWordPercentage
$ws.Cells.Item(4, 8).Value  = 45                       # H4 - negWordPercentage
$ws.Cells.Item(4, 9).Value  = 74                       # I4 - posPhrasePercentage
$ws.Cells.Item(4, 10).Value = 24                       # J4 - negPhrasePercentage
$ws.Cells.Item(4, 11).Value = 9187                     # K4 - ElapsedMs
$ws.Cells.Item(4, 12).Value = 280                      # L4 - posWordCount
$ws.Cells.Item(4, 13).Value = 251                      # M4 - negWordCount
$ws.Cells.Item(4, 14).Value = 18                       # N4 - positivePhraseCount
$ws.Cells.Item(4, 15).Value = 6                        # O4 - negativePhraseCount
$ws.Cells.Item(4, 16).Value = "Bag"                     # P4 - Method
$ws.Cells.Item(4, 17).Value = 0                        # Q4 - RSI
$ws.Cells.Item(4, 18).Value = 0.49                     # R4 - PEG
$ws.Cells.Item(4, 19).Value = 0.088800000000000004     # S4 - 200Moving%
$ws.Cells.Item(4, 19).NumberFormat = "0.00%"
$ws.Cells.Item(4, 20).Value = -1.1000000000000001      # T4 - 50Moving%
$ws.Cells.Item(4, 21).Value = 2.2999999999999998       # U4 - PriceBook
$ws.Cells.Item(4, 22).Value = "N/A"                     # V4 - Dividend
$ws.Cells.Item(4, 23).Value = 0                        # W4 - Bollinger

# --- Add new columns X and Y data to existing row 3 (Bollinger trend delta + category) ---
$ws.Cells.Item(3, 24).Value = -0.19999900000000181   # X3
$ws.Cells.Item(3, 25).Value = "Down"                  # Y3

# --- Column C width tweak (to fit the new "Neutral" text) ---
$ws.Columns.Item(3).ColumnWidth = 6.8
